$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 6 (Monday, 2021-03-22): add date, tweak Time Out + Total hours
# ---------------------------------------------------------------------
$ws.Range("B6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 44277
$ws.Range("D6").Value = 0.5
$ws.Range("E6").Value = 3

# ---------------------------------------------------------------------
# Row 7 (Tuesday, 2021-03-23): fill in the whole row
# ---------------------------------------------------------------------
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B7").Value = 44278

$ws.Range("C7").Value = 0.66666666666666663
$ws.Range("D7").Value = 0.75
$ws.Range("E7").Value = 2

$ws.Range("F6").Copy() | Out-Null
$ws.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F7").Value = "Work on the first part of pitch presentation slides"

$ws.Range("F6").Copy() | Out-Null
$ws.Range("G7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G7").Value = "Prepare materials required for pitch presentation"

$ws.Range("H6").Copy() | Out-Null
$ws.Range("H7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H7").WrapText = $true
$ws.Range("H7").Value = "Complete the first draft of the first part of the presentation slides"

# ---------------------------------------------------------------------
# Row 9 (Thursday, 2021-03-25): fill in the whole row, taller row height
# ---------------------------------------------------------------------
$ws.Rows("9").RowHeight = 50

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B9").Value = 44280

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C9").Value = 0.625

$ws.Range("D6").Copy() | Out-Null
$ws.Range("D9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D9").Value = 0.75

$ws.Range("E9").Value = 3

$ws.Range("F9").Value = "Work on the 'What' and 'Why' of the business case, first draft milestone plan "

$ws.Range("F6").Copy() | Out-Null
$ws.Range("G9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G9").Value = "As part of the business case and first milestone plan draft "

$ws.Range("H6").Copy() | Out-Null
$ws.Range("H9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H9").WrapText = $true
$ws.Range("H9").Value = "Complete the draft of most of the 'What' and 'Why' of the business case, first milestone plan"

# ---------------------------------------------------------------------
# Row 10 (Friday, 2021-03-26): add date, re-colour time cells
# ---------------------------------------------------------------------
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B10").Value = 44281

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C10").Value = 0.875

$ws.Range("D6").Copy() | Out-Null
$ws.Range("D10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D10").Value = 0.625

# ---------------------------------------------------------------------
# Selection + recalculated total
# ---------------------------------------------------------------------
$ws.Range("B9").Select() | Out-Null
$ws.Calculate() | Out-Null

Write-Output "done"
